# Generate Report for Handoff
# Updates the localization-status report: marks b.md as "Ready for handoff"
# in all three sheets, and records the new handoff xliff files/timestamps
# plus the stale-handback error detail for b.md.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-31 08:45:29"

# --- zh-cn sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "Ready for handoff"
# "Content Duplicate" must stay literal text "False" (not boolean) - copy
# from F2, which already holds the literal text "False", to avoid Excel's
# auto-boolean conversion when assigning the string directly.
$zhcn.Range("F2").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-31 08:45:24"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8274e1c7dbb6ee85c6d364bd08b69d694c445a4a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a14276d04190a98e69bf5c4a99819583f63c2f7b/e2e/b.md."
# Widen "Error Detail" (col P) to fit the new long error text - COM
# ColumnWidth uses "characters" units; 39.17 round-trips to the OOXML
# width="40" (same unit mismatch already present on this sheet's other
# width="40" columns, e.g. G/J "Latest Handoff/Handback File").
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F2").Copy($dede.Range("F3"))
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-31 08:45:29"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8274e1c7dbb6ee85c6d364bd08b69d694c445a4a/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a14276d04190a98e69bf5c4a99819583f63c2f7b/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
